$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 173
$ws.Range("F7").Value = 10236
$ws.Range("F9").Value = 143
$ws.Range("F11").Value = 725
$ws.Range("F12").Value = 4054
$ws.Range("F14").Value = 156
$ws.Range("F15").Value = 95
$ws.Range("F16").Value = 5757
$ws.Range("F18").Value = 2249
$ws.Range("F21").Value = 422
$ws.Range("F22").Value = 8637
$ws.Range("F24").Value = 107
$ws.Range("F26").Value = 2261
$ws.Range("F27").Value = 2321
$ws.Range("F28").Value = 1361
$ws.Range("F29").Value = 207
$ws.Range("F30").Value = 1898
$ws.Range("F32").Value = 45
$ws.Range("F33").Value = 310
$ws.Range("F38").Value = 23
$ws.Range("F40").Value = 1200
$ws.Range("F42").Value = 82
$ws.Range("F43").Value = 210
$ws.Range("F44").Value = 1438
$ws.Range("F45").Value = 2321
$ws.Range("F46").Value = 179
$ws.Range("F47").Value = 265
$ws.Range("F48").Value = 1235
$ws.Range("F49").Value = 19
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 162
$ws.Range("F10").Value = 933
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 653
$ws.Range("F3").Value = 848
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 653
$ws.Range("F4").Value = 848
$ws.Range("F7").Value = 10236
$ws.Range("F9").Value = 143
$ws.Range("F10").Value = 162
$ws.Range("F13").Value = 725
$ws.Range("F14").Value = 4054
$ws.Range("F15").Value = 35
$ws.Range("F16").Value = 156
$ws.Range("F17").Value = 95
$ws.Range("F18").Value = 5757
$ws.Range("F21").Value = 422
$ws.Range("F22").Value = 8637
$ws.Range("F25").Value = 107
$ws.Range("F27").Value = 2261
$ws.Range("F28").Value = 1361
$ws.Range("F29").Value = 207
$ws.Range("F30").Value = 1898
$ws.Range("F31").Value = 45
$ws.Range("F32").Value = 310
$ws.Range("F36").Value = 23
$ws.Range("F37").Value = 1200
$ws.Range("F38").Value = 82
$ws.Range("F39").Value = 210
$ws.Range("F40").Value = 1438
$ws.Range("F42").Value = 2321
$ws.Range("F43").Value = 179
$ws.Range("F45").Value = 265
$ws.Range("F48").Value = 1235
$ws.Range("F49").Value = 19
